$d = $word.ActiveDocument

$replacements = @(
    @{old = "37×38=1406"; new = "87×60=5220"},
    @{old = "11×13=143"; new = "76×56=4256"},
    @{old = "33×45=1485"; new = "47×18=846"},
    @{old = "49×86=4214"; new = "16×83=1328"},
    @{old = "86×49=4214"; new = "61×88=5368"},
    @{old = "61×28=1708"; new = "47×68=3196"},
    @{old = "42×82=3444"; new = "59×25=1475"},
    @{old = "74×82=6068"; new = "52×59=3068"},
    @{old = "93×75=6975"; new = "83×57=4731"},
    @{old = "96×34=3264"; new = "28×85=2380"},
    @{old = "93×22=2046"; new = "12×21=252"},
    @{old = "89×12=1068"; new = "59×79=4661"},
    @{old = "40×90=3600"; new = "60×18=1080"},
    @{old = "85×31=2635"; new = "60×58=3480"},
    @{old = "14×45=630"; new = "41×76=3116"},
    @{old = "18×87=1566"; new = "63×68=4284"},
    @{old = "83×13=1079"; new = "53×18=954"},
    @{old = "91×37=3367"; new = "19×23=437"},
    @{old = "98×25=2450"; new = "38×69=2622"},
    @{old = "81×63=5103"; new = "75×12=900"},
    @{old = "28×76=2128"; new = "47×53=2491"},
    @{old = "22×59=1298"; new = "91×79=7189"},
    @{old = "82×45=3690"; new = "37×62=2294"},
    @{old = "76×26=1976"; new = "74×28=2072"},
    @{old = "46×31=1426"; new = "80×89=7120"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
